$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring the header / id cell values in line with MarkdownDirectoryTransformer's
# naming convention (drop the "foaf:" prefix, use "@id" instead of "id", and
# qualify the person identifier with "ss-person:").
$ws.Range("A1").Value = "@id"
$ws.Range("B1").Value = "familyName"
$ws.Range("D1").Value = "name"
$ws.Range("A2").Value = "ss-person:minor-gordon"
